# Refresh the cryptocurrency snapshot: updates the "Price" (D) and
# "Volume(1h)" (E) columns for each coin row (2-51) with the latest scraped
# values, as produced by the GitHub Actions cron job.
#
# D/E cells in this sheet are plain text (e.g. "1.001", "  +0.20%  "), not
# numbers. Assigning a numeric-looking string straight to Range.Value makes
# Excel auto-convert it to a real number (e.g. "1.000" -> 1), so we flip the
# Price column to a Text number format before writing, then restore the
# default "Normal" style afterwards so the saved cells look just like the
# originals (no explicit style index).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range('D2').Value = '25.999.59'
$ws.Range('E2').Value = '  -0.28%  '
$ws.Range('D3').Value = '1.745.72'
$ws.Range('E3').Value = '  -0.10%  '
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '249.14'
$ws.Range('E5').Value = '  +6.43%  '
$ws.Range('E6').Value = '  +0.12%  '
$ws.Range('D7').Value = '0.5145'
$ws.Range('E7').Value = '  -1.97%  '
$ws.Range('D8').Value = '0.2756'
$ws.Range('E8').Value = '  -1.14%  '
$ws.Range('D9').Value = '0.06198'
$ws.Range('E9').Value = '  +0.01%  '
$ws.Range('D10').Value = '1.741.89'
$ws.Range('E10').Value = '  -1.03%  '
$ws.Range('E11').Value = '  +0.92%  '
$ws.Range('D12').Value = '15.18'
$ws.Range('E12').Value = '  -1.32%  '
$ws.Range('D13').Value = '0.6488'
$ws.Range('E13').Value = '  +0.55%  '
$ws.Range('D14').Value = '4.633'
$ws.Range('E14').Value = '  +1.08%  '
$ws.Range('D15').Value = '77.83'
$ws.Range('E15').Value = '  -0.50%  '
$ws.Range('E16').Value = '  +0.18%  '
$ws.Range('D17').Value = '1.000'
$ws.Range('E17').Value = '  +0.03%  '
$ws.Range('D18').Value = '26.032.11'
$ws.Range('E18').Value = '  +0.21%  '
$ws.Range('D19').Value = '11.86'
$ws.Range('E19').Value = '  +1.67%  '
$ws.Range('D20').Value = '0.000006832'
$ws.Range('E20').Value = '  +1.84%  '
$ws.Range('D21').Value = '1.964.20'
$ws.Range('E21').Value = '  -0.67%  '
$ws.Range('D22').Value = '4.295'
$ws.Range('E22').Value = '  -0.29%  '
$ws.Range('D23').Value = '8.683'
$ws.Range('E23').Value = '  -1.76%  '
$ws.Range('D24').Value = '5.366'
$ws.Range('E24').Value = '  +3.01%  '
$ws.Range('D25').Value = '135.26'
$ws.Range('E25').Value = '  -3.45%  '
$ws.Range('D26').Value = '1.503'
$ws.Range('E26').Value = '  -1.14%  '
$ws.Range('D27').Value = '15.28'
$ws.Range('E27').Value = '  -0.45%  '
$ws.Range('E28').Value = '  -1.69%  '
$ws.Range('D29').Value = '105.97'
$ws.Range('E29').Value = '  +1.70%  '
$ws.Range('D30').Value = '3.957'
$ws.Range('E30').Value = '  +4.59%  '
$ws.Range('D31').Value = '0.08257'
$ws.Range('E31').Value = '  -1.20%  '
$ws.Range('D32').Value = '3.679'
$ws.Range('E32').Value = '  -0.40%  '
$ws.Range('E33').Value = '  +2.69%  '
$ws.Range('D34').Value = '2.654'
$ws.Range('E34').Value = '  +0.71%  '
$ws.Range('D35').Value = '0.9993'
$ws.Range('E35').Value = '  -0.06%  '
$ws.Range('D36').Value = '0.6238'
$ws.Range('E36').Value = '  -1.01%  '
$ws.Range('D37').Value = '2.735'
$ws.Range('E37').Value = '  +1.01%  '
$ws.Range('D38').Value = '0.01605'
$ws.Range('E38').Value = '  +0.36%  '
$ws.Range('D39').Value = '1.933'
$ws.Range('E39').Value = '  -0.52%  '
$ws.Range('D40').Value = '1.000'
$ws.Range('E40').Value = '  +0.11%  '
$ws.Range('D41').Value = '100.46'
$ws.Range('E41').Value = '  +1.82%  '
$ws.Range('D42').Value = '0.3883'
$ws.Range('E42').Value = '  -0.60%  '
$ws.Range('D43').Value = '0.7567'
$ws.Range('E43').Value = '  +2.66%  '
$ws.Range('D44').Value = '5.022'
$ws.Range('E44').Value = '  -1.02%  '
$ws.Range('D45').Value = '6.342'
$ws.Range('E45').Value = '  +0.40%  '
$ws.Range('D46').Value = '0.1134'
$ws.Range('E46').Value = '  -0.31%  '
$ws.Range('D47').Value = '55.46'
$ws.Range('E47').Value = '  +2.79%  '
$ws.Range('D48').Value = '0.05229'
$ws.Range('E48').Value = '  -2.46%  '
$ws.Range('D49').Value = '30.70'
$ws.Range('E49').Value = '  +0.86%  '
$ws.Range('D50').Value = '7.607'
$ws.Range('E50').Value = '  -0.81%  '
$ws.Range('D51').Value = '0.3444'
$ws.Range('E51').Value = '  -0.79%  '

$priceRange.Style = "Normal"
